$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Sheet4")

# Insert a new column before the HuMFre column (K) for the new
# "Biological risk assessment number" field (x1246).
$ws4.Columns("K").Insert()

$ws4.Range("K2").Value = "Biological risk assessment number"
$ws4.Range("K3").Value = "RISKX"

# Match the header/value formatting used by the neighbouring columns.
$ws4.Range("L2").Copy()
$ws4.Range("K2").PasteSpecial(-4122)
$ws4.Range("L3").Copy()
$ws4.Range("K3").PasteSpecial(-4122)

# Sheet4 becomes the active sheet/tab, with K3 selected.
$ws4.Activate()
$ws4.Range("K3").Select()
